$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1010.2308
$ws.Range("I32").Value = 720.2
$ws.Range("J32").Value = 1191.5
$ws.Range("K32").Value = 720.2
$ws.Range("L32").Value = 1191.5
$ws.Range("M32").Value = -394.2
$ws.Range("N32").Value = -1843.5

$ws.Range("H64").Value = 3362.3584
$ws.Range("I64").Value = 3242.65
$ws.Range("J64").Value = 3730.6924
$ws.Range("K64").Value = 3242.65
$ws.Range("L64").Value = 3730.6924
$ws.Range("M64").Value = -2994.65
$ws.Range("N64").Value = -4226.6924

$ws.Range("H67").Value = 3362.3584
$ws.Range("I67").Value = 3242.65
$ws.Range("J67").Value = 3730.6924
$ws.Range("K67").Value = 3242.65
$ws.Range("L67").Value = 3730.6924
$ws.Range("M67").Value = -2384.65
$ws.Range("N67").Value = -5446.6924

$ws.Range("H129").Value = 783.9400000000001
$ws.Range("J129").Value = 811.5106
$ws.Range("L129").Value = 2434.5318
$ws.Range("N129").Value = -12434.5318

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()

$ws.Range("H58").Value = 14500
$ws.Range("J58").Value = 14500
$ws.Range("L58").Value = 14500
$ws.Range("N58").Value = -15360

$ws.Range("H59").Value = 9222.223
$ws.Range("I59").Value = 8833.333000000001
$ws.Range("K59").Value = 8833.333000000001
$ws.Range("M59").Value = -8029.333000000001

$ws.Range("H60").Value = 9320.4
$ws.Range("I60").Value = 8867.333000000001
$ws.Range("J60").Value = 10000
$ws.Range("K60").Value = 8867.333000000001
$ws.Range("L60").Value = 10000
$ws.Range("M60").Value = -8134.333000000001
$ws.Range("N60").Value = -11466

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws.Range("H76").Value = 55555
$ws.Range("J76").Value = 55555
$ws.Range("L76").Value = 55555
$ws.Range("N76").Value = -56231

$ws.Range("H79").Value = 55555
$ws.Range("J79").Value = 55555
$ws.Range("L79").Value = 55555
$ws.Range("N79").Value = -57895

$ws.Range("H80").Value = 32620
$ws.Range("J80").Value = 32620
$ws.Range("L80").Value = 32620
$ws.Range("N80").Value = -34616

$ws.Range("H83").Value = 32620
$ws.Range("J83").Value = 32620
$ws.Range("L83").Value = 97860
$ws.Range("N83").Value = -107844

$ws.Range("H87").Value = 30000
$ws.Range("J87").Value = 30000
$ws.Range("L87").Value = 30000
$ws.Range("N87").Value = -32496

$ws.Range("H88").Value = 1595
$ws.Range("I88").Value = 1456.6666
$ws.Range("J88").Value = 1733.3334
$ws.Range("K88").Value = 1456.6666
$ws.Range("L88").Value = 1733.3334
$ws.Range("M88").Value = -1050.6666
$ws.Range("N88").Value = -2545.3334

$ws.Range("H90").Value = 30000
$ws.Range("J90").Value = 30000
$ws.Range("L90").Value = 90000
$ws.Range("N90").Value = -102480

$ws.Range("H91").Value = 1595
$ws.Range("I91").Value = 1456.6666
$ws.Range("J91").Value = 1733.3334
$ws.Range("K91").Value = 1456.6666
$ws.Range("L91").Value = 1733.3334
$ws.Range("M91").Value = -52.66660000000002
$ws.Range("N91").Value = -4541.3334

$ws.Range("H92").Value = 35385
$ws.Range("J92").Value = 35385
$ws.Range("L92").Value = 35385
$ws.Range("N92").Value = -40377

$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()

$ws.Range("H95").Value = 41000
$ws.Range("J95").Value = 41000
$ws.Range("L95").Value = 41000
$ws.Range("N95").Value = -46492

$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 38666.668
$ws.Range("J68").Value = 38666.668
$ws.Range("L68").Value = 38666.668
$ws.Range("N68").Value = -40288.668

$ws.Range("H71").Value = 38666.668
$ws.Range("J71").Value = 38666.668
$ws.Range("L71").Value = 116000.004
$ws.Range("N71").Value = -124112.004

$ws.Range("H76").Value = 10000
$ws.Range("J76").Value = 10000
$ws.Range("L76").Value = 10000
$ws.Range("N76").Value = -10630

$ws.Range("H79").Value = 10000
$ws.Range("J79").Value = 10000
$ws.Range("L79").Value = 10000
$ws.Range("N79").Value = -12184

$ws.Range("H94").Value = 539.7
$ws.Range("I94").Value = 490.1
$ws.Range("J94").Value = 638.9
$ws.Range("K94").Value = 490.1
$ws.Range("L94").Value = 638.9
$ws.Range("M94").Value = -39.10000000000002
$ws.Range("N94").Value = -1540.9

$ws.Range("H99").Value = 2808.05
$ws.Range("I99").Value = 2015.5555
$ws.Range("J99").Value = 3456.4546
$ws.Range("K99").Value = 2015.5555
$ws.Range("L99").Value = 3456.4546
$ws.Range("M99").Value = -517.5554999999999
$ws.Range("N99").Value = -6452.4546

$ws.Range("H133").Value = 30640
$ws.Range("J133").Value = 30640
$ws.Range("L133").Value = 30640
$ws.Range("N133").Value = -40760

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2698.5103
$ws.Range("I31").Value = 1850.3422
$ws.Range("J31").Value = 5628.5454
$ws.Range("K31").Value = 1850.3422
$ws.Range("L31").Value = 5628.5454
$ws.Range("M31").Value = -1555.3422
$ws.Range("N31").Value = -6218.5454

$ws.Range("H34").Value = 2698.5103
$ws.Range("I34").Value = 1850.3422
$ws.Range("J34").Value = 5628.5454
$ws.Range("K34").Value = 1850.3422
$ws.Range("L34").Value = 5628.5454
$ws.Range("M34").Value = -1648.3422
$ws.Range("N34").Value = -6032.5454

$ws.Range("H99").Value = 2487.3845
$ws.Range("I99").Value = 2165.375
$ws.Range("J99").Value = 3002.6
$ws.Range("K99").Value = 2165.375
$ws.Range("L99").Value = 3002.6
$ws.Range("M99").Value = -667.375
$ws.Range("N99").Value = -5998.6

$ws.Range("H126").Value = 2487.3845
$ws.Range("I126").Value = 2165.375
$ws.Range("J126").Value = 3002.6
$ws.Range("K126").Value = 6496.125
$ws.Range("L126").Value = 9007.799999999999
$ws.Range("M126").Value = -4026.125
$ws.Range("N126").Value = -13947.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 8000
$ws.Range("I56").Value = 8000
$ws.Range("K56").Value = 8000
$ws.Range("M56").Value = -7470

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1430.8918
$ws.Range("I102").Value = 1333.8438
$ws.Range("J102").Value = 2052
$ws.Range("K102").Value = 1333.8438
$ws.Range("L102").Value = 2052
$ws.Range("M102").Value = 288.1561999999999
$ws.Range("N102").Value = -5296

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3362.5264
$ws.Range("I40").Value = 2218.8
$ws.Range("K40").Value = 2218.8
$ws.Range("M40").Value = -2082.8

$ws.Range("H136").Value = 3405.4482
$ws.Range("I136").Value = 2011.125
$ws.Range("J136").Value = 5121.5386
$ws.Range("K136").Value = 6033.375
$ws.Range("L136").Value = 15364.6158
$ws.Range("M136").Value = -3483.375
$ws.Range("N136").Value = -20464.6158

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3320.2827
$ws.Range("I136").Value = 3304.3057
$ws.Range("J136").Value = 3377.8
$ws.Range("K136").Value = 9912.917099999999
$ws.Range("L136").Value = 10133.4
$ws.Range("M136").Value = -7362.917099999999
$ws.Range("N136").Value = -15233.4
